$wb = $excel.ActiveWorkbook

# --- Rename the "Include" worksheets ---
$wsEntity = $wb.Worksheets.Item("Include from EntityCode")
$wsEntity.Name = "Include #0"

$wsSpecimen = $wb.Worksheets.Item("Include from SpecimenType")
$wsSpecimen.Name = "Include #1"

# --- Update the Metadata sheet ---
$ws1 = $wb.Worksheets.Item("Metadata")

# Update Version value
$ws1.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Update Date value
$ws1.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new row for "Jurisdiction" right before the "Description" row,
# matching the formatting of the surrounding data rows.
$ws1.Rows.Item(11).Insert()
$ws1.Range("A12:B12").Copy()
$ws1.Range("A11:B11").PasteSpecial(-4122)
$ws1.Range("A11").Value = "Jurisdiction"

# B11 must stay a (shared) empty string value, same as the existing blank
# "Value" cells elsewhere in the workbook, rather than a truly blank cell.
$wsEntity.Range("B104").Copy()
$ws1.Range("B11").PasteSpecial(-4163)

Write-Output "done"
